$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4318.4
$ws.Range("I40").Value = 2435
$ws.Range("J40").Value = 5574
$ws.Range("K40").Value = 2435
$ws.Range("L40").Value = 5574
$ws.Range("M40").Value = -2260
$ws.Range("N40").Value = -5924

$ws.Range("H129").Value = 1091.2
$ws.Range("I129").Value = 406.5
$ws.Range("J129").Value = 1262.375
$ws.Range("K129").Value = 1219.5
$ws.Range("L129").Value = 3787.125
$ws.Range("M129").Value = 3780.5
$ws.Range("N129").Value = -13787.125

$ws.Range("H137").Value = 1898.3334
$ws.Range("I137").Value = 1168.619
$ws.Range("J137").Value = 3601
$ws.Range("K137").Value = 3505.857
$ws.Range("L137").Value = 10803
$ws.Range("M137").Value = -955.857
$ws.Range("N137").Value = -15903

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 1800
$ws.Range("I29").Value = 1800
$ws.Range("K29").Value = 1800
$ws.Range("M29").Value = -1492

$ws.Range("H32").Value = 8697.053
$ws.Range("I32").Value = 2399.9421
$ws.Range("J32").Value = 25408.615
$ws.Range("K32").Value = 2399.9421
$ws.Range("L32").Value = 25408.615
$ws.Range("M32").Value = -2112.9421
$ws.Range("N32").Value = -25982.615

$ws.Range("H74").Value = 11367697
$ws.Range("I74").Value = 22732212
$ws.Range("J74").Value = 3182.3635
$ws.Range("K74").Value = 22732212
$ws.Range("L74").Value = 3182.3635
$ws.Range("M74").Value = -22731338
$ws.Range("N74").Value = -4930.363499999999

$ws.Range("H77").Value = 11367697
$ws.Range("I77").Value = 22732212
$ws.Range("J77").Value = 3182.3635
$ws.Range("K77").Value = 113661060
$ws.Range("L77").Value = 15911.8175
$ws.Range("M77").Value = -113656692
$ws.Range("N77").Value = -24647.8175

$ws.Range("H132").Value = 1505.3877
$ws.Range("I132").Value = 1230.0555
$ws.Range("J132").Value = 2267.8462
$ws.Range("K132").Value = 3690.1665
$ws.Range("L132").Value = 6803.5386
$ws.Range("M132").Value = -1160.1665
$ws.Range("N132").Value = -11863.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 78028.39999999999
$ws.Range("J9").Value = 78028.39999999999
$ws.Range("L9").Value = 78028.39999999999
$ws.Range("N9").Value = -78364.39999999999

$ws.Range("H31").Value = 3655779
$ws.Range("I31").Value = 5290407.5
$ws.Range("J31").Value = 1903.1765
$ws.Range("K31").Value = 5290407.5
$ws.Range("L31").Value = 1903.1765
$ws.Range("M31").Value = -5290112.5
$ws.Range("N31").Value = -2493.1765

$ws.Range("H34").Value = 3655779
$ws.Range("I34").Value = 5290407.5
$ws.Range("J34").Value = 1903.1765
$ws.Range("K34").Value = 5290407.5
$ws.Range("L34").Value = 1903.1765
$ws.Range("M34").Value = -5290205.5
$ws.Range("N34").Value = -2307.1765

$ws.Range("H58").Value = 1219.4736
$ws.Range("I58").Value = 701.5714
$ws.Range("J58").Value = 1859.2354
$ws.Range("K58").Value = 701.5714
$ws.Range("L58").Value = 1859.2354
$ws.Range("M58").Value = -498.5714
$ws.Range("N58").Value = -2265.2354

$ws.Range("H136").Value = 1219.4736
$ws.Range("I136").Value = 701.5714
$ws.Range("J136").Value = 1859.2354
$ws.Range("K136").Value = 2104.7142
$ws.Range("L136").Value = 5577.706200000001
$ws.Range("M136").Value = 445.2857999999997
$ws.Range("N136").Value = -10677.7062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 192
$ws.Range("I107").Value = 192.66667
$ws.Range("J107").Value = 190
$ws.Range("K107").Value = 578.00001
$ws.Range("L107").Value = 570
$ws.Range("M107").Value = 1341.99999
$ws.Range("N107").Value = -4410

$ws.Range("H113").Value = 714.5493
$ws.Range("I113").Value = 559.78125
$ws.Range("J113").Value = 841.53845
$ws.Range("K113").Value = 1679.34375
$ws.Range("L113").Value = 2524.61535
$ws.Range("M113").Value = 490.65625
$ws.Range("N113").Value = -6864.61535

$ws.Range("H131").Value = 8772726
$ws.Range("I131").Value = 315
$ws.Range("J131").Value = 11495198
$ws.Range("K131").Value = 945
$ws.Range("L131").Value = 34485594
$ws.Range("M131").Value = 4095
$ws.Range("N131").Value = -34495674

$ws.Range("H132").Value = 1112.5
$ws.Range("I132").Value = 820
$ws.Range("K132").Value = 7380
$ws.Range("M132").Value = -4850

$ws.Range("H137").Value = 6179.9614
$ws.Range("I137").Value = 850
$ws.Range("J137").Value = 6624.125
$ws.Range("K137").Value = 2550
$ws.Range("L137").Value = 19872.375
$ws.Range("M137").Value = 2550
$ws.Range("N137").Value = -30072.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2540.83
$ws.Range("I132").Value = 2296.7317
$ws.Range("J132").Value = 3374.8333
$ws.Range("K132").Value = 6890.195099999999
$ws.Range("L132").Value = 10124.4999
$ws.Range("M132").Value = -4360.195099999999
$ws.Range("N132").Value = -15184.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5051887
$ws.Range("I82").Value = 7576764.5
$ws.Range("J82").Value = 2132.1667
$ws.Range("K82").Value = 7576764.5
$ws.Range("L82").Value = 2132.1667
$ws.Range("M82").Value = -7576403.5
$ws.Range("N82").Value = -2854.1667

$ws.Range("H85").Value = 5051887
$ws.Range("I85").Value = 7576764.5
$ws.Range("J85").Value = 2132.1667
$ws.Range("K85").Value = 7576764.5
$ws.Range("L85").Value = 2132.1667
$ws.Range("M85").Value = -7575516.5
$ws.Range("N85").Value = -4628.1667

$ws.Range("H93").Value = 1020.381
$ws.Range("I93").Value = 1087.5
$ws.Range("K93").Value = 1087.5
$ws.Range("M93").Value = 160.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 10170
$ws.Range("I29").Value = 490
$ws.Range("J29").Value = 19850
$ws.Range("K29").Value = 490
$ws.Range("L29").Value = 19850
$ws.Range("M29").Value = -200
$ws.Range("N29").Value = -20430

$ws.Range("H68").Value = 28757
$ws.Range("J68").Value = 28757
$ws.Range("L68").Value = 28757
$ws.Range("N68").Value = -30379

$ws.Range("H71").Value = 28757
$ws.Range("J71").Value = 28757
$ws.Range("L71").Value = 86271
$ws.Range("N71").Value = -94383

$ws.Range("H81").Value = 62502890
$ws.Range("I81").Value = 200002860
$ws.Range("J81").Value = 2899.6365
$ws.Range("K81").Value = 400005720
$ws.Range("L81").Value = 5799.273
$ws.Range("M81").Value = -400004659
$ws.Range("N81").Value = -7921.273

$ws.Range("H84").Value = 62502890
$ws.Range("I84").Value = 200002860
$ws.Range("J84").Value = 2899.6365
$ws.Range("K84").Value = 2000028600
$ws.Range("L84").Value = 28996.365
$ws.Range("M84").Value = -2000023296
$ws.Range("N84").Value = -39604.36500000001

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H132").Value = 1708.9395
$ws.Range("I132").Value = 1067.92
$ws.Range("J132").Value = 3712.125
$ws.Range("K132").Value = 3203.76
$ws.Range("L132").Value = 11136.375
$ws.Range("M132").Value = -673.7600000000002
$ws.Range("N132").Value = -16196.375
